$wb = $excel.ActiveWorkbook

# --- Rename the existing (only) sheet to "ConflictMatrix" ---
$conflictMatrix = $wb.Worksheets.Item(1)
$conflictMatrix.Name = "ConflictMatrix"

# Select A1:A4 on ConflictMatrix before creating/activating the new sheet,
# so its sheetView ends up with sqref="A1:A4" (and is not the "active" tab).
$conflictMatrix.Range("A1:A4").Select() | Out-Null

# --- Insert a brand new "Setup" sheet before ConflictMatrix ---
$setup = $wb.Worksheets.Add($conflictMatrix)
$setup.Name = "Setup"

# --- Fill the Setup sheet with its data ---
$setup.Range("A1").Value = 2.1
$setup.Range("B1").Value = 1.5

$setup.Range("A2").Value = 5.1
$setup.Range("B2").Value = 1.5

$setup.Range("A3").Value = 8.1
$setup.Range("B3").Value = 1.5

$setup.Range("A4").Value = 11.1
$setup.Range("B4").Value = 1.5

# Give column A on Setup the same "highlighted" look used elsewhere in the
# workbook: one decimal place, thin border, yellow fill.
$highlight = $setup.Range("A1:A4")
$highlight.NumberFormat = '_-* #,##0.0_-;\-* #,##0.0_-;_-* "-"??_-;_-@_-'
$highlight.Interior.Color = 10284031
$highlight.Borders.LineStyle = 1

Write-Output "Workbook now has sheets: $($wb.Worksheets.Item(1).Name), $($wb.Worksheets.Item(2).Name)"
